$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text so values like
# "332.40" or "0.4720" keep their exact original formatting/type
# instead of being auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.682.57'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '1.871.45'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '332.40'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').Value = '0.4720'
$ws.Range('E7').Value = '  +4.39%  '
$ws.Range('D8').Value = '0.3943'
$ws.Range('D9').Value = '47.93'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('D11').Value = '1.026'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '22.03'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('D13').Value = '1.870.28'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = '5.954'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '7.151'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '0.00001047'
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('D18').Value = '86.84'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').Value = '0.06633'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').Value = '17.22'
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').Value = '27.684.42'
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('D23').Value = '5.506'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = '10.98'
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('D25').Value = '2.308'
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('D26').Value = '2.091.13'
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('D27').Value = '158.93'
$ws.Range('E27').Value = '  +3.92%  '
$ws.Range('D28').Value = '20.25'
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D29').Value = '2.097'
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('D30').Value = '5.567'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = '122.32'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').Value = '0.9707'
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('D33').Value = '0.09517'
$ws.Range('E33').Value = '  +2.23%  '
$ws.Range('D34').Value = '1.451'
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').Value = '3.589'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = '5.335'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '0.06108'
$ws.Range('E37').Value = '  +1.87%  '
$ws.Range('D38').Value = '0.02257'
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('D39').Value = '1.232'
$ws.Range('D40').Value = '8.161'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('D41').Value = '0.6021'
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').Value = '0.1901'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('D44').Value = '10.28'
$ws.Range('E44').Value = '  +1.33%  '
$ws.Range('D45').Value = '1.257'
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('D46').Value = '0.5720'
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('D47').Value = '12.24'
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('D48').Value = '1.943'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').Value = '3.387'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').Value = '0.06856'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').Value = '114.65'
$ws.Range('E51').Value = '  +6.24%  '
